$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as plain text in this sheet (e.g. "66.573.36"
# uses dots as thousands separators, and small numbers like "561.73" are literal
# strings, not numeric cells). Cells whose new value would otherwise be auto-detected
# as a number are pre-formatted as Text so Excel keeps them as strings.

$ws.Range("D2").Value = "66.573.36"
$ws.Range("E2").Value = "  +5.42%  "

$ws.Range("D3").Value = "3.509.67"
$ws.Range("E3").Value = "  +8.63%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.73"
$ws.Range("E5").Value = "  +8.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.57"
$ws.Range("E6").Value = "  +10.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +10.48%  "

$ws.Range("D8").Value = "3.505.43"
$ws.Range("E8").Value = "  +8.75%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.642"
$ws.Range("E10").Value = "  +9.31%  "

$ws.Range("E11").Value = "  +20.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.80"
$ws.Range("E12").Value = "  +8.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("E13").Value = "  +11.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.42"
$ws.Range("E14").Value = "  +7.88%  "

$ws.Range("D15").Value = "4.073.70"
$ws.Range("E15").Value = "  +9.92%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.511.37"
$ws.Range("E16").Value = "  +9.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.68"
$ws.Range("E17").Value = "  +10.21%  "

$ws.Range("E18").Value = "  +6.62%  "

$ws.Range("D19").Value = "66.573.22"
$ws.Range("E19").Value = "  +6.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("E20").Value = "  +11.00%  "

$ws.Range("E21").Value = "  +7.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "420.38"
$ws.Range("E22").Value = "  +14.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +14.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.10"
$ws.Range("E24").Value = "  +7.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.17"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.01"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("E27").Value = "  +11.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.37"
$ws.Range("E28").Value = "  +13.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.11"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.18"
$ws.Range("E30").Value = "  +16.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.47"
$ws.Range("E31").Value = "  +9.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("E32").Value = "  +4.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "621.17"
$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.88"
$ws.Range("E34").Value = "  +8.63%  "

$ws.Range("E35").Value = "  +9.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.98"
$ws.Range("E36").Value = "  +6.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("E37").Value = "  +24.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "38.43"
$ws.Range("E38").Value = "  +10.30%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0813"
$ws.Range("E39").Value = "  +16.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  +5.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  +12.08%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.132.34"
$ws.Range("E43").Value = "  +11.68%  "

$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.64"
$ws.Range("E45").Value = "  +2.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  +13.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0418"
$ws.Range("E47").Value = "  +8.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  +11.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.73"
$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.134"
$ws.Range("E50").Value = "  +9.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.10"
$ws.Range("E51").Value = "  +3.50%  "
